$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, shifting existing rows 260..368 down to 261..369
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with a copy of the (pre-shift) row 260 data,
# but with an updated Fecha (date) value of 44636.
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44636
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100114014
$ws.Cells.Item(260, 7).Value = "Betarraga"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 40
$ws.Cells.Item(260, 11).Value = 8000
$ws.Cells.Item(260, 12).Value = 8000
$ws.Cells.Item(260, 13).Value = 8000
$ws.Cells.Item(260, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(260, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(260, 16).Value = 667
$ws.Cells.Item(260, 17).Value = 12
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same custom date/time number format as the rest of column D
$ws.Cells.Item(260, 4).NumberFormat = $ws.Cells.Item(261, 4).NumberFormat
$ws.Cells.Item(260, 4).Value = 44636
